# Insert a new record row at row 165 in the "Hortaliza, Femacal de La Calera - Haba"
# sheet, shifting the existing rows 165-176 down to 166-177, and populate the
# new row with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 165 (existing rows 165:176 shift down to 166:177)
$ws.Rows("165:165").Insert()

# Populate the newly inserted row 165 with the new data record
$ws.Range("A165").Value = 3
$ws.Range("B165").Value = "Femacal de La Calera"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44783
$ws.Range("E165").Value = 5
$ws.Range("F165").Value = 100112026
$ws.Range("G165").Value = "Haba"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 95
$ws.Range("K165").Value = 15000
$ws.Range("L165").Value = 16000
$ws.Range("M165").Value = 15526
$ws.Range("N165").Value = "`$/saco 25 kilos"
$ws.Range("O165").Value = "Provincia de Limarí"
$ws.Range("P165").Value = 621
$ws.Range("Q165").Value = 25
$ws.Range("R165").Value = "Hortaliza"
